# The workbook tracks an "approved systems" list that gets republished
# periodically under a sheet named for the publish date, with a workbook-level
# defined name ("Sygehus_EPJ_systemer") pointing at that sheet. This edit
# rolls the publish date forward from 02-12-2025 to 05-12-2025.
#
# Renaming the sheet via the Excel object model automatically repoints the
# defined name's sheet-qualified reference (Excel keeps defined names in
# sync with sheet renames), matching both the <sheet> and <definedName>
# updates in the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Opdateret d. 05-12-2025"
